{"js": "// Replace the division-expression text in each table cell with its\n// updated counterpart, per the commit's regenerated problem set.\n// Mapping is old-text -> new-text; every old value is unique in the\n// document, so a plain body-wide search/replace for each pair is safe.\nconst replacements = [\n  [\"349\u00f74=\", \"748\u00f73=\"],\n  [\"708\u00f75=\", \"486\u00f79=\"],\n  [\"635\u00f76=\", \"215\u00f73=\"],\n  [\"913\u00f77=\", \"766\u00f76=\"],\n  [\"414\u00f74=\", \"600\u00f76=\"],\n  [\"827\u00f78=\", \"573\u00f79=\"],\n  [\"417\u00f76=\", \"309\u00f74=\"],\n  [\"983\u00f77=\", \"344\u00f79=\"],\n  [\"226\u00f77=\", \"649\u00f74=\"],\n  [\"152\u00f76=\", \"501\u00f78=\"],\n  [\"627\u00f78=\", \"397\u00f72=\"],\n  [\"958\u00f77=\", \"487\u00f78=\"],\n  [\"150\u00f77=\", \"411\u00f78=\"],\n  [\"837\u00f75=\", \"415\u00f75=\"],\n  [\"439\u00f74=\", \"195\u00f74=\"],\n  [\"576\u00f78=\", \"881\u00f73=\"],\n  [\"413\u00f75=\", \"786\u00f75=\"],\n  [\"940\u00f79=\", \"510\u00f76=\"],\n  [\"103\u00f79=\", \"975\u00f75=\"],\n  [\"928\u00f75=\", \"285\u00f74=\"],\n  [\"477\u00f76=\", \"137\u00f72=\"],\n  [\"855\u00f76=\", \"121\u00f72=\"],\n  [\"275\u00f74=\", \"586\u00f72=\"],\n  [\"737\u00f75=\", \"557\u00f72=\"],\n  [\"759\u00f79=\", \"283\u00f75=\"],\n];\n\nconst body = context.document.body;\nconst searchResults = [];\n\nfor (const [oldText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  searchResults.push(found);\n}\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const found = searchResults[i];\n  for (const range of found.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the division-expression text in each table cell with its\n# updated counterpart, per the commit's regenerated problem set.\n# Mapping is old-text -> new-text; every old value is unique in the\n# document, so a plain Find/Replace over the whole document body is\n# safe for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"349\u00f74=\", \"748\u00f73=\"),\n    @(\"708\u00f75=\", \"486\u00f79=\"),\n    @(\"635\u00f76=\", \"215\u00f73=\"),\n    @(\"913\u00f77=\", \"766\u00f76=\"),\n    @(\"414\u00f74=\", \"600\u00f76=\"),\n    @(\"827\u00f78=\", \"573\u00f79=\"),\n    @(\"417\u00f76=\", \"309\u00f74=\"),\n    @(\"983\u00f77=\", \"344\u00f79=\"),\n    @(\"226\u00f77=\", \"649\u00f74=\"),\n    @(\"152\u00f76=\", \"501\u00f78=\"),\n    @(\"627\u00f78=\", \"397\u00f72=\"),\n    @(\"958\u00f77=\", \"487\u00f78=\"),\n    @(\"150\u00f77=\", \"411\u00f78=\"),\n    @(\"837\u00f75=\", \"415\u00f75=\"),\n    @(\"439\u00f74=\", \"195\u00f74=\"),\n    @(\"576\u00f78=\", \"881\u00f73=\"),\n    @(\"413\u00f75=\", \"786\u00f75=\"),\n    @(\"940\u00f79=\", \"510\u00f76=\"),\n    @(\"103\u00f79=\", \"975\u00f75=\"),\n    @(\"928\u00f75=\", \"285\u00f74=\"),\n    @(\"477\u00f76=\", \"137\u00f72=\"),\n    @(\"855\u00f76=\", \"121\u00f72=\"),\n    @(\"275\u00f74=\", \"586\u00f72=\"),\n    @(\"737\u00f75=\", \"557\u00f72=\"),\n    @(\"759\u00f79=\", \"283\u00f75=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
